# Applies the commit:
#  - rename first sheet "06020b" -> "blad1" and empty out its contents
#  - rewrite / extend the data on sheet "06020" (second sheet)

$wb = $excel.ActiveWorkbook

# --- 1. "06020b" -> "blad1", emptied out -------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "blad1"
$ws1.Cells.ClearContents()

# --- 2. "06020" sheet gets new / changed rows ---------------------------
$ws2 = $wb.Worksheets.Item(2)

function Set-Cell($ws, $addr, $value) {
    $ws.Range($addr).Value = $value
}

function Set-TextCell($ws, $addr, $value) {
    # Force text storage, even for numeric-/currency-looking strings, so
    # that e.g. "10142" or "€ 10043" are kept as text rather than being
    # auto-converted into a number by Excel.
    $ws.Range($addr).Value = "'" + $value
}

# Row 2 - "Soort verzekering"
Set-Cell     $ws2 "A2" "Soort verzekering"
Set-TextCell $ws2 "B2" "10142"
Set-Cell     $ws2 "F2" "02 Soort verzekering          10142"
Set-Cell     $ws2 "G2" "x"

# Row 3
Set-TextCell $ws2 "B3" "10142"
Set-Cell     $ws2 "C3" "Omschrijving"
Set-Cell     $ws2 "D3" "Links"
Set-Cell     $ws2 "E3" "niet verwijderen"

# Row 4 - "Gezinssamenstelling" (used to live on sheet1)
Set-Cell     $ws2 "A4" "Gezinssamenstelling"
Set-TextCell $ws2 "B4" "10694"
$ws2.Range("C4").ClearContents()
$ws2.Range("D4").ClearContents()
$ws2.Range("E4").ClearContents()
Set-Cell     $ws2 "F4" "03 Gezinssamenstelling        10694"
Set-Cell     $ws2 "G4" "x"

# Row 5
Set-TextCell $ws2 "B5" "10694"
Set-Cell     $ws2 "C5" "Omschrijving"
Set-Cell     $ws2 "D5" "Links"
Set-Cell     $ws2 "E5" "verwijderen"

# Row 6 - "Meeverzekerd sterdekking"
Set-Cell     $ws2 "A6" "Meeverzekerd sterdekking"
Set-TextCell $ws2 "B6" "11500"
Set-Cell     $ws2 "F6" "04 Meeverzekerd sterdekking   11500"
Set-Cell     $ws2 "G6" "x"

# Row 7
Set-TextCell $ws2 "B7" "11500"
Set-Cell     $ws2 "D7" "Links"
Set-Cell     $ws2 "E7" "verwijderen"

# Row 8 - "Meeverzekerd verhaalsbijstand"
Set-Cell     $ws2 "A8" "Meeverzekerd verhaalsbijstand"
Set-TextCell $ws2 "B8" "11646"
Set-Cell     $ws2 "F8" "05 Meeverzekerd verhaalsbijstand 11646"
Set-Cell     $ws2 "G8" "x"

# Row 9
Set-TextCell $ws2 "B9" "11646"
Set-Cell     $ws2 "D9" "Links"
Set-Cell     $ws2 "E9" "verwijderen"

# Row 10 - "Eigen risico"
Set-Cell     $ws2 "A10" "Eigen risico"
Set-TextCell $ws2 "B10" "€ 10043"
Set-Cell     $ws2 "F10" "06 Eigen risico               € 10043"
Set-Cell     $ws2 "G10" "x"

# Row 11
Set-TextCell $ws2 "B11" "10043"
Set-Cell     $ws2 "C11" "Getal inclusief decimalen"
Set-Cell     $ws2 "D11" "Links"
Set-Cell     $ws2 "E11" "verwijderen"

# Row 12 - "Eigen risico kinderen/huisdieren"
Set-Cell     $ws2 "A12" "Eigen risico kinderen/huisdieren"
Set-TextCell $ws2 "B12" "€ 11450"
Set-Cell     $ws2 "F12" "07 Eigen risico kinderen/huisdieren€ 11450"
Set-Cell     $ws2 "G12" "x"

# Row 13
Set-TextCell $ws2 "B13" "11450"
Set-Cell     $ws2 "C13" "Getal inclusief decimalen"
Set-Cell     $ws2 "D13" "Links"
Set-Cell     $ws2 "E13" "verwijderen"

# Row 14 - "Verzekerde rubrieken"
Set-Cell     $ws2 "A14" "Verzekerde rubrieken"
Set-Cell     $ws2 "B14" "                        Verzekerd bedrag"
Set-Cell     $ws2 "F14" "09 Verzekerde rubrieken                             Verzekerd bedrag"
Set-Cell     $ws2 "G14" "x"
